$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# REPORTDATE
$ws.Range("H2").Value = "2020-09-30 00:00:00"

# BASIC_EPS
$ws.Range("I2").Value = 0.54

# TOTAL_OPERATE_INCOME
$ws.Range("K2").Value = 99369292.31999999

# PARENT_NETPROFIT
$ws.Range("L2").Value = 16094268.55

# YSTZ, SJLTZ, BPS, MGJYXJJE (were blank, now numeric)
$ws.Range("N2").Value = 122.9209529084
$ws.Range("O2").Value = 348.5532237172
$ws.Range("P2").Value = 3.307587800667
$ws.Range("Q2").Value = -0.393075698667

# XSMLL
$ws.Range("R2").Value = 30.2504255874

# ISNEW - keep as text "1" (was text "0")
$ws.Range("AB2").NumberFormat = "@"
$ws.Range("AB2").Value = "1"
$ws.Range("AB2").Style = "Normal"

# QDATE
$ws.Range("AC2").Value = "2020Q3"

# DATATYPE
$ws.Range("AD2").Value = "2020年 三季报"

# DATAYEAR - keep as text "2020" (was text "2019")
$ws.Range("AE2").NumberFormat = "@"
$ws.Range("AE2").Value = "2020"
$ws.Range("AE2").Style = "Normal"
